# Update Work Week and Social Spending
# GDP per Capita data refresh for Senegal: revised historical series (1950-2010)
# plus newly reported years 2011-2016.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Column E holds the GDP per Capita figures as text (matches the source
# workbook, where these numeric-looking values are stored as strings).
$ws.Range("E2:E68").NumberFormat = "@"

# row, year, revised GDP per Capita value
$gdpRows = @(
    @(2, 1950, "2007"),
    @(3, 1951, "2042"),
    @(4, 1952, "2077"),
    @(5, 1953, "2110"),
    @(6, 1954, "2146"),
    @(7, 1955, "2179"),
    @(8, 1956, "2212"),
    @(9, 1957, "2246"),
    @(10, 1958, "2278"),
    @(11, 1959, "2308"),
    @(12, 1960, "2303"),
    @(13, 1961, "2351"),
    @(14, 1962, "2370"),
    @(15, 1963, "2402"),
    @(16, 1964, "2389"),
    @(17, 1965, "2409"),
    @(18, 1966, "2404"),
    @(19, 1967, "2310"),
    @(20, 1968, "2389"),
    @(21, 1969, "2171"),
    @(22, 1970, "2287"),
    @(23, 1971, "2216"),
    @(24, 1972, "2289"),
    @(25, 1973, "2096"),
    @(26, 1974, "2118"),
    @(27, 1975, "2225"),
    @(28, 1976, "2370"),
    @(29, 1977, "2236"),
    @(30, 1978, "2082"),
    @(31, 1979, "2158"),
    @(32, 1980, "2021"),
    @(33, 1981, "1943"),
    @(34, 1982, "2168"),
    @(35, 1983, "2153"),
    @(36, 1984, "1989"),
    @(37, 1985, "2000"),
    @(38, 1986, "2031"),
    @(39, 1987, "2048"),
    @(40, 1988, "2085"),
    @(41, 1989, "1980"),
    @(42, 1990, "1999"),
    @(43, 1991, "1982.62595787466"),
    @(44, 1992, "1942.92672653003"),
    @(45, 1993, "1906.82386892074"),
    @(46, 1994, "1847.59827373069"),
    @(47, 1995, "1886.54281966775"),
    @(48, 1996, "1866.34228502269"),
    @(49, 1997, "1866.39106308534"),
    @(50, 1998, "1915.42965790462"),
    @(51, 1999, "1978.09642159315"),
    @(52, 2000, "1983.52404748961"),
    @(53, 2001, "2009.94228319219"),
    @(54, 2002, "1959.61409887651"),
    @(55, 2003, "2025.02345473418"),
    @(56, 2004, "2074.95381162079"),
    @(57, 2005, "2121.1586016431"),
    @(58, 2006, "2104.43690696656"),
    @(59, 2007, "2137.67164678246"),
    @(60, 2008, "2146.79146464314"),
    @(61, 2009, "2129.33639483785"),
    @(62, 2010, "2151.32470542")
)

foreach ($entry in $gdpRows) {
    $row = $entry[0]
    $ws.Cells.Item($row, 5).Value = $entry[2]
}

# Newly added years 2011-2016
$newRows = @(
    @(63, 2011, "2123"),
    @(64, 2012, "2163"),
    @(65, 2013, "2185"),
    @(66, 2014, "2224"),
    @(67, 2015, "2311"),
    @(68, 2016, "2403")
)

foreach ($entry in $newRows) {
    $row = $entry[0]
    $year = $entry[1]
    $val = $entry[2]
    $ws.Cells.Item($row, 1).Value = 686
    $ws.Cells.Item($row, 2).Value = "Senegal"
    $ws.Cells.Item($row, 3).Value = "GDP per Capita"
    $ws.Cells.Item($row, 4).Value = $year
    $ws.Cells.Item($row, 5).Value = $val
}
